# Generate Report for handback
# The source file "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.md" has finished its
# handback round-trip (zh-cn and de-de). Update the localization-status
# report: flip its Status from "Ready for handoff" to
# "Handed back: in sync with en-US", and fill in the Latest Target File /
# Latest Handback File / Latest Handback DateTime columns on the per-locale
# sheets.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: row 5 is the e38ef9a3 file; both locale columns (B, C)
# move from "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B5").Value = $statusHandedBack
$wsOverview.Range("C5").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet: row 5 is the e38ef9a3 file.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B5").Value = $statusHandedBack
$wsZh.Range("E5").Value = "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.md"
$wsZh.Range("F5").Value = "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.247e64e92cc46294a4588ea0f04c2e248e6b4e58.zh-cn.xlf"
$wsZh.Range("G5").Value = "2016-01-22 02:39:38"

$wsZh.Hyperlinks.Add($wsZh.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0aac78fd0909119d71b277128a808df3a99c4fd8/e2e/e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.md", "", "", "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/9b99be4f4b430084bdd5c289c4fc501f02c96114/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.247e64e92cc46294a4588ea0f04c2e248e6b4e58.zh-cn.xlf", "", "", "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.247e64e92cc46294a4588ea0f04c2e248e6b4e58.zh-cn.xlf")

$wsZh.Range("E5").Font.Underline = 2
$wsZh.Range("E5").Font.Color = 15570276
$wsZh.Range("F5").Font.Underline = 2
$wsZh.Range("F5").Font.Color = 15570276

# ---------------------------------------------------------------------
# de-de sheet: row 5 is the e38ef9a3 file.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B5").Value = $statusHandedBack
$wsDe.Range("E5").Value = "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.md"
$wsDe.Range("F5").Value = "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.247e64e92cc46294a4588ea0f04c2e248e6b4e58.de-de.xlf"
$wsDe.Range("G5").Value = "2016-01-22 02:40:02"

$wsDe.Hyperlinks.Add($wsDe.Range("E5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bf17789921a763334c0f94f487670a907017accd/e2e/e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.md", "", "", "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b0ff5fa5172feb8b02c92cda4cd1bc1927509237/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.247e64e92cc46294a4588ea0f04c2e248e6b4e58.de-de.xlf", "", "", "e38ef9a3-98b8-4466-ba5a-4c7a06a844ef.247e64e92cc46294a4588ea0f04c2e248e6b4e58.de-de.xlf")

$wsDe.Range("E5").Font.Underline = 2
$wsDe.Range("E5").Font.Color = 15570276
$wsDe.Range("F5").Font.Underline = 2
$wsDe.Range("F5").Font.Color = 15570276
